# "Conversa sobre a heuristica 2"
#
# Clears out the per-algorithm statistics (node counts, timings, branching
# factor, penetrance, etc.) from the "Tabuleiro A", "Tabuleiro B" and
# "Tabuleiro C" sheets, leaving only the row labels in column A (and, on
# "Tabuleiro C", the "indefinido" note in B2). The previously-applied cell
# formats (date format on F4, left-aligned G2/G4 on "Tabuleiro B") are kept
# even though the values are gone, since those cells keep their styling.

$wb = $excel.ActiveWorkbook

# --- "Tabuleiro A": wipe all stats except the algorithm names in column A ---
$wsA = $wb.Worksheets.Item("Tabuleiro A")
$wsA.Range("B2:G5").ClearContents() | Out-Null
$wsA.Range("H3").ClearContents() | Out-Null
$wsA.Range("B2").Select() | Out-Null

# --- "Tabuleiro B": same cleanup, selection becomes the whole data block ---
$wsB = $wb.Worksheets.Item("Tabuleiro B")
$wsB.Range("B2:G5").ClearContents() | Out-Null
$wsB.Range("H3").ClearContents() | Out-Null
$wsB.Range("B2:H5").Select() | Out-Null

# --- "Tabuleiro C": wipe stats but keep B2's "indefinido" note ---
$wsC = $wb.Worksheets.Item("Tabuleiro C")
$wsC.Range("B3:G5").ClearContents() | Out-Null
$wsC.Range("H3").ClearContents() | Out-Null
$wsC.Range("G4").Select() | Out-Null

# "Tabuleiro D", "Tabuleiro E" and "Tabuleiro F" already only held the
# algorithm labels, so nothing changes there.

# Re-select "Tabuleiro A" so it stays the visible/active sheet.
$wsA.Activate() | Out-Null
